$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3482.3555
$ws.Range("I15").Value = 3482.3555
$ws.Range("K15").Value = 10447.0665
$ws.Range("M15").Value = -10278.0665

$ws.Range("H17").Value = 2128595
$ws.Range("J17").Value = 2174855.8
$ws.Range("L17").Value = 6524567.399999999
$ws.Range("N17").Value = -6524903.399999999

$ws.Range("H28").Value = 499.5
$ws.Range("I28").Value = 499.5
$ws.Range("K28").Value = 499.5
$ws.Range("M28").Value = -14.5

$ws.Range("H33").Value = 980.63635
$ws.Range("I33").Value = 754.1111
$ws.Range("J33").Value = 2000
$ws.Range("K33").Value = 754.1111
$ws.Range("L33").Value = 2000
$ws.Range("M33").Value = -525.1111
$ws.Range("N33").Value = -2458

$ws.Range("H43").Value = 2314
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2314
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2314
$ws.Range("N43").Value = -2452
$ws.Range("M43").ClearContents()

$ws.Range("H64").Value = 10841.429
$ws.Range("I64").Value = 2945
$ws.Range("J64").Value = 14000
$ws.Range("K64").Value = 2945
$ws.Range("L64").Value = 14000
$ws.Range("M64").Value = -2697
$ws.Range("N64").Value = -14496

$ws.Range("H67").Value = 10841.429
$ws.Range("I67").Value = 2945
$ws.Range("J67").Value = 14000
$ws.Range("K67").Value = 2945
$ws.Range("L67").Value = 14000
$ws.Range("M67").Value = -2087
$ws.Range("N67").Value = -15716

$ws.Range("H112").Value = 3622.111
$ws.Range("J112").Value = 4142.857
$ws.Range("L112").Value = 12428.571
$ws.Range("N112").Value = -14644.571

$ws.Range("H138").Value = 6519.7095
$ws.Range("J138").Value = 6550.95
$ws.Range("L138").Value = 19652.85
$ws.Range("N138").Value = -29932.85

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2628.3845
$ws.Range("J2").Value = 2768.5
$ws.Range("L2").Value = 2768.5
$ws.Range("N2").Value = -2994.5

$ws.Range("H45").Value = 290404.44
$ws.Range("I45").Value = 668166
$ws.Range("J45").Value = 7083.25
$ws.Range("K45").Value = 668166
$ws.Range("L45").Value = 7083.25
$ws.Range("M45").Value = -667789
$ws.Range("N45").Value = -7837.25

$ws.Range("H61").Value = 3270605.8
$ws.Range("I61").Value = 3706014
$ws.Range("J61").Value = 5043.6665
$ws.Range("K61").Value = 3706014
$ws.Range("L61").Value = 5043.6665
$ws.Range("M61").Value = -3705802
$ws.Range("N61").Value = -5467.6665

$ws.Range("H116").Value = 2628.3845
$ws.Range("J116").Value = 2768.5
$ws.Range("L116").Value = 2768.5
$ws.Range("N116").Value = -7356.5

$ws.Range("H136").Value = 3270605.8
$ws.Range("I136").Value = 3706014
$ws.Range("J136").Value = 5043.6665
$ws.Range("K136").Value = 11118042
$ws.Range("L136").Value = 15130.9995
$ws.Range("M136").Value = -11115492
$ws.Range("N136").Value = -20230.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2628.3845
$ws.Range("J3").Value = 2768.5
$ws.Range("L3").Value = 2768.5
$ws.Range("N3").Value = -2996.5

$ws.Range("H134").Value = 3152.1614
$ws.Range("J134").Value = 2814.5
$ws.Range("L134").Value = 8443.5
$ws.Range("N134").Value = -13513.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 32116.666
$ws.Range("J28").Value = 32116.666
$ws.Range("L28").Value = 32116.666
$ws.Range("N28").Value = -32606.666

$ws.Range("H31").Value = 5415.8335
$ws.Range("J31").Value = 6566.857
$ws.Range("L31").Value = 6566.857
$ws.Range("N31").Value = -7156.857

$ws.Range("H34").Value = 5415.8335
$ws.Range("J34").Value = 6566.857
$ws.Range("L34").Value = 6566.857
$ws.Range("N34").Value = -6970.857

$ws.Range("H74").Value = 42603
$ws.Range("J74").Value = 42603
$ws.Range("L74").Value = 42603
$ws.Range("N74").Value = -44351

$ws.Range("H77").Value = 42603
$ws.Range("J77").Value = 42603
$ws.Range("L77").Value = 127809
$ws.Range("N77").Value = -136545

$ws.Range("H81").Value = 62000
$ws.Range("J81").Value = 62000
$ws.Range("L81").Value = 62000
$ws.Range("N81").Value = -63996

$ws.Range("H84").Value = 62000
$ws.Range("J84").Value = 62000
$ws.Range("L84").Value = 186000
$ws.Range("N84").Value = -195984

$ws.Range("H97").Value = 27195.8
$ws.Range("J97").Value = 26494.75
$ws.Range("L97").Value = 26494.75
$ws.Range("N97").Value = -28476.75

$ws.Range("H111").Value = 72912.5
$ws.Range("J111").Value = 72912.5
$ws.Range("L111").Value = 72912.5
$ws.Range("N111").Value = -81092.5

$ws.Range("H112").Value = 64917.25
$ws.Range("J112").Value = 64917.25
$ws.Range("L112").Value = 64917.25
$ws.Range("N112").Value = -67871.25

$ws.Range("H134").Value = 9229.5
$ws.Range("J134").Value = 11699.286
$ws.Range("L134").Value = 35097.858
$ws.Range("N134").Value = -40167.858

$ws.Range("H141").Value = 47766.5
$ws.Range("J141").Value = 47766.5
$ws.Range("L141").Value = 47766.5
$ws.Range("N141").Value = -58126.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 13159968
$ws.Range("I68").Value = 20835144
$ws.Range("J68").Value = 2522.6428
$ws.Range("K68").Value = 62505432
$ws.Range("L68").Value = 7567.928400000001
$ws.Range("M68").Value = -62504621
$ws.Range("N68").Value = -9189.928400000001

$ws.Range("H71").Value = 13159968
$ws.Range("I71").Value = 20835144
$ws.Range("J71").Value = 2522.6428
$ws.Range("K71").Value = 187516296
$ws.Range("L71").Value = 22703.7852
$ws.Range("M71").Value = -187512240
$ws.Range("N71").Value = -30815.7852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 17955
$ws.Range("J99").Value = 22000
$ws.Range("L99").Value = 22000
$ws.Range("N99").Value = -26492

$ws.Range("H107").Value = 887.4
$ws.Range("I107").Value = 887.4
$ws.Range("K107").Value = 887.4
$ws.Range("M107").Value = 1032.6

$ws.Range("H122").Value = 3300.2
$ws.Range("I122").Value = 3300.2
$ws.Range("K122").Value = 9900.599999999999
$ws.Range("M122").Value = -7450.599999999999

$ws.Range("H132").Value = 2457.7273
$ws.Range("I132").Value = 2457.7273
$ws.Range("K132").Value = 7373.1819
$ws.Range("M132").Value = -4843.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2818.2273
$ws.Range("I22").Value = 2080.5454
$ws.Range("K22").Value = 2080.5454
$ws.Range("M22").Value = -1785.5454

$ws.Range("H27").Value = 2818.2273
$ws.Range("I27").Value = 2080.5454
$ws.Range("K27").Value = 2080.5454
$ws.Range("M27").Value = -1973.5454

$ws.Range("H34").Value = 22380.715
$ws.Range("J34").Value = 15000
$ws.Range("L34").Value = 15000
$ws.Range("N34").Value = -15344

$ws.Range("H46").Value = 2430.2307
$ws.Range("I46").Value = 1519.6
$ws.Range("K46").Value = 1519.6
$ws.Range("M46").Value = -1331.6

$ws.Range("H95").Value = 39088.6
$ws.Range("J95").Value = 39088.6
$ws.Range("L95").Value = 39088.6
$ws.Range("N95").Value = -44580.6

$ws.Range("H136").Value = 4780.647
$ws.Range("I136").Value = 3523.0833
$ws.Range("J136").Value = 7798.8
$ws.Range("K136").Value = 10569.2499
$ws.Range("L136").Value = 23396.4
$ws.Range("M136").Value = -8019.249899999999
$ws.Range("N136").Value = -28496.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 14393.637
$ws.Range("J45").Value = 14393.637
$ws.Range("L45").Value = 14393.637
$ws.Range("N45").Value = -15375.637

$ws.Range("H54").Value = 27405
$ws.Range("J54").Value = 27405
$ws.Range("L54").Value = 27405
$ws.Range("N54").Value = -28445
